$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114..147 down to 115..148.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with a new data record
# (same fixed dims as the rest of the sheet; new Fecha/volume/price data).
$ws.Cells.Item(114, 1).Value = 7
$ws.Cells.Item(114, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(114, 3).Value = 'Ñuble'
$ws.Cells.Item(114, 4).Value = 44463
$ws.Cells.Item(114, 5).Value = 16
$ws.Cells.Item(114, 6).Value = 100112043
$ws.Cells.Item(114, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(114, 8).Value = 'Sin especificar'
$ws.Cells.Item(114, 9).Value = 'Primera'
$ws.Cells.Item(114, 10).Value = 300
$ws.Cells.Item(114, 11).Value = 16000
$ws.Cells.Item(114, 12).Value = 17000
$ws.Cells.Item(114, 13).Value = 16500
$ws.Cells.Item(114, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(114, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(114, 16).Value = 275
$ws.Cells.Item(114, 17).Value = 60
$ws.Cells.Item(114, 18).Value = 'Hortaliza'

# Ensure the sheet dimension reflects the new extent.
$ws.Range("A1:R148").Select()
